# Normalize the "Recorded By" column (G): ensure the literal "System"
# entry (exact case) always appears first in the comma-separated list.
# For the one special case without a "System" entry
# ("admin@admin.com, dnasr281@gmail.com") swap the order instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    $val = [string]$val
    if ($val -eq "") {
        continue
    }

    $parts = @($val -split ", ")

    $hasExactSystem = $false
    foreach ($p in $parts) {
        if ($p.Equals("System")) {
            $hasExactSystem = $true
        }
    }

    if ($hasExactSystem) {
        $newParts = @()
        $newParts += "System"
        $alreadyRemoved = $false
        foreach ($p in $parts) {
            if ((-not $alreadyRemoved) -and $p.Equals("System")) {
                $alreadyRemoved = $true
            }
            else {
                $newParts += $p
            }
        }
        $newVal = $newParts -join ", "
        if (-not $newVal.Equals($val)) {
            $cell.Value2 = $newVal
        }
    }
    elseif ($val.Equals("admin@admin.com, dnasr281@gmail.com")) {
        $cell.Value2 = "dnasr281@gmail.com, admin@admin.com"
    }
}
